$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing "SmartRules" example block (D16:E20) into a new
# block at D31:E35, representing the second/overloaded rules signature
# that shares the same token aliases (EPBDS-11473).
$src = $ws.Range("D16:E20")
$dst = $ws.Range("D31:E35")
$src.Copy($dst)

# Give the new block its own header text (second overload signature).
$ws.Range("D31").Value = "SmartRules MyDatatype myRules2( MyDatatype myObj)"

# Merge the header row across D31:E31, same as the look of similar
# example headers in this workbook.
$null = $ws.Range("D31:E31").Merge()

# Match the saved selection/active cell shown in the workbook.
$null = $ws.Range("D31:E31").Select()
